$wb = $excel.ActiveWorkbook

# Update the "Status" value from "Ready for handoff" to "In Translation"
# across the Overview sheet (columns E and F) and the per-language sheets (column C).

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$overview.Range("E2:F2").ColumnWidth = 13.4101845877511

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"
$zhcn.Range("C2").EntireColumn.ColumnWidth = 13.4101845877511

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"
$dede.Range("C2").EntireColumn.ColumnWidth = 13.4101845877511
